$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.251.71'
$ws.Range("E2").Value = '  -0.16%  '
$ws.Range("D3").Value = '1.592.22'
$ws.Range("E3").Value = '  +0.18%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''212.87'
$ws.Range("E5").Value = '  +0.40%  '
$ws.Range("E6").Value = '  -0.42%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -0.44%  '
$ws.Range("D9").Value = '''0.0606'
$ws.Range("E9").Value = '  -0.47%  '
$ws.Range("D10").Value = '''18.96'
$ws.Range("E10").Value = '  -2.06%  '
$ws.Range("D11").Value = '''0.0849'
$ws.Range("E11").Value = '  +0.24%  '
$ws.Range("D12").Value = '1.818.30'
$ws.Range("E12").Value = '  +0.29%  '
$ws.Range("D13").Value = '1.590.42'
$ws.Range("E13").Value = '  -2.75%  '
$ws.Range("D14").Value = '''3.99'
$ws.Range("E14").Value = '  -1.12%  '
$ws.Range("E15").Value = '  -2.29%  '
$ws.Range("D16").Value = '''63.77'
$ws.Range("E16").Value = '  -0.98%  '
$ws.Range("D17").Value = '26.257.92'
$ws.Range("E17").Value = '  -0.20%  '
$ws.Range("D18").Value = '0.0₃0723'
$ws.Range("E18").Value = '  -1.12%  '
$ws.Range("D19").Value = '''215.29'
$ws.Range("E19").Value = '  +0.78%  '
$ws.Range("E20").Value = '  -2.00%  '
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("E22").Value = '  +0.23%  '
$ws.Range("D23").Value = '''9.02'
$ws.Range("E23").Value = '  +0.41%  '
$ws.Range("E24").Value = '  -2.16%  '
$ws.Range("D25").Value = '''144.81'
$ws.Range("E25").Value = '  -0.21%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").Value = '''6.95'
$ws.Range("E27").Value = '  -1.24%  '
$ws.Range("E28").Value = '  +0.85%  '
$ws.Range("D29").Value = '''15.09'
$ws.Range("E29").Value = '  -0.55%  '
$ws.Range("E30").Value = '  -1.72%  '
$ws.Range("D31").Value = '''1.15'
$ws.Range("E31").Value = '  -0.03%  '
$ws.Range("E32").Value = '  -0.45%  '
$ws.Range("D33").Value = '1.419.18'
$ws.Range("E33").Value = '  +5.81%  '
$ws.Range("E35").Value = '  -0.89%  '
$ws.Range("E36").Value = '  -1.64%  '
$ws.Range("D37").Value = '''0.575'
$ws.Range("E37").Value = '  -4.06%  '
$ws.Range("E38").Value = '  -0.87%  '
$ws.Range("E39").Value = '  +0.83%  '
$ws.Range("D40").Value = '''5.78'
$ws.Range("E40").Value = '  +0.11%  '
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("B42").Value = 'WEMIXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D42").Value = '''0.939'
$ws.Range("E42").Value = '  -11.13%  '
$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").Value = '''2.16'
$ws.Range("E43").Value = '  +0.78%  '
$ws.Range("D44").Value = '''0.760'
$ws.Range("E44").Value = '  -0.24%  '
$ws.Range("D45").Value = '1.730.75'
$ws.Range("E45").Value = '  +0.32%  '
$ws.Range("D46").Value = '''60.85'
$ws.Range("E46").Value = '  -1.38%  '
$ws.Range("D47").Value = '''86.86'
$ws.Range("E47").Value = '  -1.45%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0104'
$ws.Range("E48").Value = '  -0.76%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = '''1.48'
$ws.Range("E49").Value = '  -1.11%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '''0.0501'
$ws.Range("E50").Value = '  -0.31%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").Value = '''0.0951'
$ws.Range("E51").Value = '  -2.94%  '
